# Auto-generated edit script: update Jenova Profits sheets per diff
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2998.4
$ws.Range("I19").Value = 999
$ws.Range("J19").Value = 3498.25
$ws.Range("K19").Value = 999
$ws.Range("L19").Value = 3498.25
$ws.Range("M19").Value = -824
$ws.Range("N19").Value = -3848.25
$ws.Range("H38").Value = 2729.25
$ws.Range("I38").Value = 1661.5625
$ws.Range("K38").Value = 4984.6875
$ws.Range("M38").Value = -4612.6875
$ws.Range("H42").Value = 309.30768
$ws.Range("I42").Value = 122.375
$ws.Range("J42").Value = 608.4
$ws.Range("K42").Value = 367.125
$ws.Range("L42").Value = 1825.2
$ws.Range("M42").Value = -137.125
$ws.Range("N42").Value = -2285.2
$ws.Range("H58").Value = 4016.8572
$ws.Range("I58").Value = 2714
$ws.Range("J58").Value = 5201.273
$ws.Range("K58").Value = 8142
$ws.Range("L58").Value = 15603.819
$ws.Range("M58").Value = -7992
$ws.Range("N58").Value = -15903.819
$ws.Range("H128").Value = 89653.44500000001
$ws.Range("J128").Value = 89653.44500000001
$ws.Range("L128").Value = 89653.44500000001
$ws.Range("N128").Value = -99613.44500000001
$ws.Range("H131").Value = 4494.8423
$ws.Range("I131").Value = 2879.4443
$ws.Range("J131").Value = 5948.7
$ws.Range("K131").Value = 8638.332900000001
$ws.Range("L131").Value = 17846.1
$ws.Range("M131").Value = -3598.332900000001
$ws.Range("N131").Value = -27926.1
$ws.Range("H137").Value = 3233.0544
$ws.Range("I137").Value = 1821.3125
$ws.Range("K137").Value = 5463.9375
$ws.Range("M137").Value = -2913.9375
$ws.Range("H138").Value = 4398.55
$ws.Range("J138").Value = 4857.164
$ws.Range("L138").Value = 14571.492
$ws.Range("N138").Value = -24851.492
$ws.Range("H141").Value = 2003.8572
$ws.Range("I141").Value = 2003.8572
$ws.Range("K141").Value = 6011.571599999999
$ws.Range("M141").Value = -831.5715999999993

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 2489.6667
$ws.Range("I3").Value = 2500
$ws.Range("J3").Value = 2484.5
$ws.Range("K3").Value = 2500
$ws.Range("L3").Value = 2484.5
$ws.Range("M3").Value = -2385
$ws.Range("N3").Value = -2714.5
$ws.Range("H32").Value = 2703.7407
$ws.Range("I32").Value = 2301.8
$ws.Range("K32").Value = 2301.8
$ws.Range("M32").Value = -2014.8
$ws.Range("H132").Value = 6332.75
$ws.Range("I132").Value = 2156.1177
$ws.Range("J132").Value = 11066.267
$ws.Range("K132").Value = 6468.353099999999
$ws.Range("L132").Value = 33198.801
$ws.Range("M132").Value = -3938.353099999999
$ws.Range("N132").Value = -38258.801

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1262.1666
$ws.Range("J64").Value = 1288.7778
$ws.Range("L64").Value = 1288.7778
$ws.Range("N64").Value = -1738.7778
$ws.Range("H67").Value = 1262.1666
$ws.Range("J67").Value = 1288.7778
$ws.Range("L67").Value = 1288.7778
$ws.Range("N67").Value = -2848.7778
$ws.Range("H86").Value = 2835968.5
$ws.Range("I86").Value = 5668766.5
$ws.Range("J86").Value = 3170.3333
$ws.Range("K86").Value = 5668766.5
$ws.Range("L86").Value = 3170.3333
$ws.Range("M86").Value = -5667643.5
$ws.Range("N86").Value = -5416.3333
$ws.Range("H89").Value = 2835968.5
$ws.Range("I89").Value = 5668766.5
$ws.Range("J89").Value = 3170.3333
$ws.Range("K89").Value = 28343832.5
$ws.Range("L89").Value = 15851.6665
$ws.Range("M89").Value = -28338216.5
$ws.Range("N89").Value = -27083.6665
$ws.Range("H107").Value = 2021.8889
$ws.Range("I107").Value = 2049.75
$ws.Range("K107").Value = 2049.75
$ws.Range("M107").Value = -129.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 7155.8237
$ws.Range("I16").Value = 2286.2727
$ws.Range("J16").Value = 16083.333
$ws.Range("K16").Value = 2286.2727
$ws.Range("L16").Value = 16083.333
$ws.Range("M16").Value = -1999.2727
$ws.Range("N16").Value = -16657.333
$ws.Range("H31").Value = 2506877.2
$ws.Range("I31").Value = 10000012
$ws.Range("J31").Value = 9165.666999999999
$ws.Range("K31").Value = 10000012
$ws.Range("L31").Value = 9165.666999999999
$ws.Range("M31").Value = -9999717
$ws.Range("N31").Value = -9755.666999999999
$ws.Range("H34").Value = 2506877.2
$ws.Range("I34").Value = 10000012
$ws.Range("J34").Value = 9165.666999999999
$ws.Range("K34").Value = 10000012
$ws.Range("L34").Value = 9165.666999999999
$ws.Range("M34").Value = -9999810
$ws.Range("N34").Value = -9569.666999999999
$ws.Range("H47").Value = 34355
$ws.Range("I47").Value = 50000
$ws.Range("J47").Value = 29140
$ws.Range("K47").Value = 50000
$ws.Range("L47").Value = 29140
$ws.Range("M47").Value = -49434
$ws.Range("N47").Value = -30272
$ws.Range("H109").Value = 74989.5
$ws.Range("J109").Value = 74989.5
$ws.Range("L109").Value = 74989.5
$ws.Range("N109").Value = -77069.5
$ws.Range("H113").Value = 7155.8237
$ws.Range("I113").Value = 2286.2727
$ws.Range("J113").Value = 16083.333
$ws.Range("K113").Value = 2286.2727
$ws.Range("L113").Value = 16083.333
$ws.Range("M113").Value = -116.2727
$ws.Range("N113").Value = -20423.333

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1260
$ws.Range("I98").Value = 1729.5
$ws.Range("K98").Value = 5188.5
$ws.Range("M98").Value = -3690.5
$ws.Range("H129").Value = 6201.273
$ws.Range("J129").Value = 9338.357
$ws.Range("L129").Value = 28015.071
$ws.Range("N129").Value = -38015.071

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 250005300
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").Value = ""
$ws.Range("H73").Value = 250005300
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").Value = ""
$ws.Range("H113").Value = 1259975.9
$ws.Range("I113").Value = 1437116
$ws.Range("K113").Value = 1437116
$ws.Range("M113").Value = -1434946
$ws.Range("H122").Value = 2937.7856
$ws.Range("J122").Value = 6500
$ws.Range("L122").Value = 19500
$ws.Range("N122").Value = -24400
$ws.Range("H132").Value = 482236.78
$ws.Range("I132").Value = 559052.6
$ws.Range("J132").Value = 205699.8
$ws.Range("K132").Value = 1677157.8
$ws.Range("L132").Value = 617099.3999999999
$ws.Range("M132").Value = -1674627.8
$ws.Range("N132").Value = -622159.3999999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").Value = ""
$ws.Range("H61").Value = 5222.8076
$ws.Range("I61").Value = 4199.8423
$ws.Range("K61").Value = 4199.8423
$ws.Range("M61").Value = -3997.8423
$ws.Range("H113").Value = 5222.8076
$ws.Range("I113").Value = 4199.8423
$ws.Range("K113").Value = 4199.8423
$ws.Range("M113").Value = -2029.8423
$ws.Range("H122").Value = 1679684.1
$ws.Range("I122").Value = 1446314.9
$ws.Range("K122").Value = 4338944.699999999
$ws.Range("M122").Value = -4336494.699999999
$ws.Range("H123").Value = 79973.5
$ws.Range("J123").Value = 79973.5
$ws.Range("L123").Value = 79973.5
$ws.Range("N123").Value = -89773.5
$ws.Range("H128").Value = 75000
$ws.Range("J128").Value = 75000
$ws.Range("L128").Value = 75000
$ws.Range("N128").Value = -84960
$ws.Range("H132").Value = 4485.36
$ws.Range("I132").Value = 3398.9333
$ws.Range("J132").Value = 6115
$ws.Range("K132").Value = 10196.7999
$ws.Range("L132").Value = 18345
$ws.Range("M132").Value = -7666.7999
$ws.Range("N132").Value = -23405

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 24896.977
$ws.Range("I132").Value = 957.65515
$ws.Range("J132").Value = 71179.664
$ws.Range("K132").Value = 2872.96545
$ws.Range("L132").Value = 213538.992
$ws.Range("M132").Value = -342.9654500000001
$ws.Range("N132").Value = -218598.992
$ws.Range("H136").Value = 68282.164
$ws.Range("I136").Value = 1289.6957
$ws.Range("K136").Value = 3869.0871
$ws.Range("M136").Value = -1319.0871

